$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF (58) holds the "Date" values (header text "Date" in BF1).
# Rows 2-31 contain the literal text "4-18-2013-14"; this was one day off
# from the correct date because of how NBA stats were shown, so replace it
# with the text "2014-04-18".
#
# Assigning a plain string like "2014-04-18" straight to .Value makes Excel
# "smart-detect" it as a date and reformat the cell, which would add a new
# number format / cell style that isn't part of the intended change. Instead
# build the text with a formula (so it's never date-parsed), then convert it
# to a plain value with copy / paste-special so the result is stored as a
# literal string without touching the cell's style.
$rng = $ws.Range($ws.Cells.Item(2, 58), $ws.Cells.Item(31, 58))
$rng.Formula = '="2014-04-18"'
$rng.Copy()
$rng.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
